# Adds a "time_taken" column (F) with per-row timestamps, mirroring the
# header styling already used on row 1 for the other columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - same style as the other header cells (B1:E1)
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 10:51:50.012209",
    "2021-10-05 10:51:50.012221",
    "2021-10-05 10:51:50.012226",
    "2021-10-05 10:51:50.012229",
    "2021-10-05 10:51:50.012232",
    "2021-10-05 10:51:50.012235",
    "2021-10-05 10:51:50.012239",
    "2021-10-05 10:51:50.012242",
    "2021-10-05 10:51:50.012245",
    "2021-10-05 10:51:50.012248",
    "2021-10-05 10:51:50.012251",
    "2021-10-05 10:51:50.012254",
    "2021-10-05 10:51:50.012257",
    "2021-10-05 10:51:50.012260",
    "2021-10-05 10:51:50.012263",
    "2021-10-05 10:51:50.012266",
    "2021-10-05 10:51:50.012269",
    "2021-10-05 10:51:50.012273",
    "2021-10-05 10:51:50.012276",
    "2021-10-05 10:51:50.012278"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
